# Auto-generated Excel COM-interop script applying the BRVM recommendations refresh
$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd  = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": refreshed stats (sector rows 2-8, stock rows 12-34) ---
$wsReco.Cells.Item(2, 3).Value = 5
$wsReco.Cells.Item(2, 4).Value = 810.73
$wsReco.Cells.Item(2, 5).Value = 164.65

$wsReco.Cells.Item(3, 3).Value = 5
$wsReco.Cells.Item(3, 4).Value = 727.83
$wsReco.Cells.Item(3, 5).Value = 145.01

$wsReco.Cells.Item(4, 3).Value = 5
$wsReco.Cells.Item(4, 4).Value = 712.78
$wsReco.Cells.Item(4, 5).Value = 142.04

$wsReco.Cells.Item(5, 3).Value = 5
$wsReco.Cells.Item(5, 4).Value = 666.3200000000001
$wsReco.Cells.Item(5, 5).Value = 136.25

$wsReco.Cells.Item(6, 3).Value = 5
$wsReco.Cells.Item(6, 4).Value = 556.6
$wsReco.Cells.Item(6, 5).Value = 112.34

$wsReco.Cells.Item(7, 3).Value = 5
$wsReco.Cells.Item(7, 4).Value = 534.55
$wsReco.Cells.Item(7, 5).Value = 107.16

$wsReco.Cells.Item(8, 3).Value = 5
$wsReco.Cells.Item(8, 4).Value = 465.62
$wsReco.Cells.Item(8, 5).Value = 93.28

$wsReco.Cells.Item(12, 2).Value = 5
$wsReco.Cells.Item(12, 4).Value = 35.82
$wsReco.Cells.Item(12, 5).Value = 7.07

$wsReco.Cells.Item(13, 1).Value = "SICOR CI (SICC)"
$wsReco.Cells.Item(13, 4).Value = 11.48
$wsReco.Cells.Item(13, 5).Value = 4.08

$wsReco.Cells.Item(14, 1).Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Cells.Item(14, 2).Value = 2
$wsReco.Cells.Item(14, 4).Value = 9.69
$wsReco.Cells.Item(14, 5).Value = 4.73

$wsReco.Cells.Item(15, 1).Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Cells.Item(15, 4).Value = 6.04
$wsReco.Cells.Item(15, 5).Value = 6.04

$wsReco.Cells.Item(16, 1).Value = "NEI-CEDA CI (NEIC)"
$wsReco.Cells.Item(16, 2).Value = 2
$wsReco.Cells.Item(16, 3).Value = 2
$wsReco.Cells.Item(16, 4).Value = 3.82
$wsReco.Cells.Item(16, 5).Value = 6.48
$wsReco.Cells.Item(16, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(19, 1).Value = "SICABLE CI (CABC)"
$wsReco.Cells.Item(19, 2).Value = 1
$wsReco.Cells.Item(19, 4).Value = 1.48
$wsReco.Cells.Item(19, 5).Value = 3.82

$wsReco.Cells.Item(22, 1).Value = "UNILEVER CI (UNLC)"
$wsReco.Cells.Item(22, 2).Value = 1
$wsReco.Cells.Item(22, 4).Value = 0.47
$wsReco.Cells.Item(22, 5).Value = -6.67
$wsReco.Cells.Item(22, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(23, 1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$wsReco.Cells.Item(23, 2).Value = 2
$wsReco.Cells.Item(23, 3).Value = 2
$wsReco.Cells.Item(23, 4).Value = 0.16
$wsReco.Cells.Item(23, 5).Value = 3.75
$wsReco.Cells.Item(23, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(24, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$wsReco.Cells.Item(24, 4).Value = -1.49
$wsReco.Cells.Item(24, 5).Value = -1.49

$wsReco.Cells.Item(25, 1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$wsReco.Cells.Item(25, 4).Value = -1.8
$wsReco.Cells.Item(25, 5).Value = -1.8

$wsReco.Cells.Item(26, 1).Value = "UNIWAX CI (UNXC)"
$wsReco.Cells.Item(26, 4).Value = -2.13
$wsReco.Cells.Item(26, 5).Value = -2.13

$wsReco.Cells.Item(27, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$wsReco.Cells.Item(27, 2).Value = 0
$wsReco.Cells.Item(27, 3).Value = 1
$wsReco.Cells.Item(27, 4).Value = -2.21
$wsReco.Cells.Item(27, 5).Value = -2.21
$wsReco.Cells.Item(27, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(28, 1).Value = "SOLIBRA CI (SLBC)"
$wsReco.Cells.Item(28, 2).Value = 2
$wsReco.Cells.Item(28, 3).Value = 2
$wsReco.Cells.Item(28, 4).Value = -2.87
$wsReco.Cells.Item(28, 5).Value = -4.33
$wsReco.Cells.Item(28, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(30, 1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$wsReco.Cells.Item(30, 3).Value = 1
$wsReco.Cells.Item(30, 4).Value = -3.92
$wsReco.Cells.Item(30, 5).Value = 3.19

$wsReco.Cells.Item(31, 1).Value = "CFAO MOTORS CI (CFAC)"
$wsReco.Cells.Item(31, 3).Value = 2
$wsReco.Cells.Item(31, 4).Value = -5.95
$wsReco.Cells.Item(31, 5).Value = 4.71

$wsReco.Cells.Item(32, 1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$wsReco.Cells.Item(32, 3).Value = 2
$wsReco.Cells.Item(32, 4).Value = -6.57
$wsReco.Cells.Item(32, 5).Value = -3.51

$wsReco.Cells.Item(33, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$wsReco.Cells.Item(33, 3).Value = 2
$wsReco.Cells.Item(33, 4).Value = -8.699999999999999
$wsReco.Cells.Item(33, 5).Value = -4.35

$wsReco.Cells.Item(34, 1).Value = "SETAO CI (STAC)"
$wsReco.Cells.Item(34, 2).Value = 0
$wsReco.Cells.Item(34, 4).Value = -12.54
$wsReco.Cells.Item(34, 5).Value = -7.05
$wsReco.Cells.Item(34, 7).Value = "➖ Neutre"

# --- Sheet "Top_YTD": refreshed YTD progression (rows 2-8) ---
$wsYtd.Cells.Item(2, 2).Value = 12277.58
$wsYtd.Cells.Item(3, 2).Value = 8829.620000000001
$wsYtd.Cells.Item(4, 2).Value = 8295.52
$wsYtd.Cells.Item(5, 2).Value = 6804.56
$wsYtd.Cells.Item(6, 2).Value = 4113.7
$wsYtd.Cells.Item(7, 2).Value = 3692.28
$wsYtd.Cells.Item(8, 2).Value = 2586.41
